$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 (leg ROM header values) updated
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 (CON) meanEMG values updated
$ws.Range("B2").Value = 107.15940310918681
$ws.Range("C2").Value = 105.0175368920161
$ws.Range("D2").Value = 105.94536616336812
$ws.Range("E2").Value = 107.15940310918681

# Row 3 (STR) meanEMG values updated; C3 cleared, D3 newly populated
$ws.Range("B3").Value = 106.66637623656203
$ws.Range("C3").ClearContents()
$ws.Range("D3").Value = 105.4767460211276
$ws.Range("E3").Value = 106.18326688907649

# Reflect the selection used while editing this range
$ws.Range("B1:E3").Select()
